# Update countries & provincias Spain
# - Refreshes the "Datos actualizados" timestamp.
# - Refreshes COVID case counters for a handful of countries with new data
#   from the source feed.
# - Because several countries' totals changed, a few rows swap position in
#   the (case-count sorted) table; since column A holds the country name and
#   columns B:H hold the stats, a swap is expressed as updating the country
#   name text in place while the stats for the row follow the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 10:35"

# --- Rusia (row 7): refreshed totals ---
$ws.Range("B7").Value = 1167805
$ws.Range("C7").Value = 8232
$ws.Range("D7").Value = 952399
$ws.Range("E7").Value = 194861
$ws.Range("G7").Value = 160
$ws.Range("H7").Value = 20545

# --- Filipinas (row 24): refreshed totals ---
$ws.Range("B24").Value = 309303
$ws.Range("C24").Value = 2025
$ws.Range("D24").Value = 252930
$ws.Range("E24").Value = 50925
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 5448

# --- Indonesia (row 26): refreshed totals ---
$ws.Range("B26").Value = 282724
$ws.Range("C26").Value = 4002
$ws.Range("D26").Value = 210437
$ws.Range("E26").Value = 61686
$ws.Range("G26").Value = 128
$ws.Range("H26").Value = 10601

# --- Polonia (row 46): refreshed totals ---
$ws.Range("D46").Value = 68955
$ws.Range("E46").Value = 17234

# --- Singapur (row 59): refreshed totals ---
$ws.Range("B59").Value = 57742
$ws.Range("C59").Value = 27
$ws.Range("E59").Value = 322

# --- Croacia overtakes Madagascar: rows 88/89 swap ---
$ws.Cells.Item(88, 1).Value = "Croacia"
$ws.Range("B88").Value = 16380
$ws.Range("C88").Value = 135
$ws.Range("E88").Value = 1158
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 275

$ws.Cells.Item(89, 1).Value = "Madagascar"
$ws.Range("B89").Value = 16348
$ws.Range("D89").Value = 14947
$ws.Range("E89").Value = 1172
$ws.Range("H89").Value = 229

# --- Eslovaquia (row 106): refreshed totals ---
$ws.Range("B106").Value = 9574
$ws.Range("C106").Value = 231
$ws.Range("D106").Value = 4329
$ws.Range("E106").Value = 5200
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 45

# --- Georgia (row 116): refreshed totals ---
$ws.Range("E116").Value = 3508
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 34

# --- Eslovenia overtakes Cuba/Suazilandia/Republica de Yibuti: rows 119-122 shift ---
$ws.Cells.Item(119, 1).Value = "Eslovenia"
$ws.Range("B119").Value = 5487
$ws.Range("C119").Value = 99
$ws.Range("D119").Value = 3682
$ws.Range("E119").Value = 1656
$ws.Range("H119").Value = 149

$ws.Cells.Item(120, 1).Value = "Cuba"
$ws.Range("B120").Value = 5483
$ws.Range("D120").Value = 4787
$ws.Range("E120").Value = 574
$ws.Range("H120").Value = 122

$ws.Cells.Item(121, 1).Value = "Suazilandia"
$ws.Range("B121").Value = 5452
$ws.Range("D121").Value = 4844
$ws.Range("E121").Value = 500
$ws.Range("H121").Value = 108

$ws.Cells.Item(122, 1).Value = "Republica de Yibuti"
$ws.Range("B122").Value = 5410
$ws.Range("D122").Value = 5340
$ws.Range("E122").Value = 9
$ws.Range("H122").Value = 61

# --- Lituania (row 131): refreshed totals ---
$ws.Range("B131").Value = 4578
$ws.Range("C131").Value = 88
$ws.Range("D131").Value = 2349
$ws.Range("E131").Value = 2137

# --- Sri Lanka (row 142): refreshed totals ---
$ws.Range("D142").Value = 3230
$ws.Range("E142").Value = 120

# --- Estonia (row 143): refreshed totals ---
$ws.Range("B143").Value = 3315
$ws.Range("C143").Value = 49
$ws.Range("D143").Value = 2564
$ws.Range("E143").Value = 687

# --- Santa Lucia / Timor Oriental tie-break swap: rows 207/208 (no stat change) ---
$ws.Cells.Item(207, 1).Value = "Santa Lucia"
$ws.Cells.Item(208, 1).Value = "Timor Oriental"
